$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 -- New Mexico
$ws.Range("B13").Value = 44041
$ws.Range("C13").Value = 20136
$ws.Range("D13").Value = 632
$ws.Range("E13").Value = 374
$ws.Range("G13").Value = 1.86

# Row 27 -- Colorado
$ws.Range("B27").Value = 44041
$ws.Range("C27").Value = 45796
$ws.Range("D27").Value = 1822
$ws.Range("G27").Value = 6.02
$ws.Range("H27").Value = 6.87
$ws.Range("K27").Value = 35396
$ws.Range("L27").Value = 1761

# Row 37 -- New Hampshire
$ws.Range("B37").Value = 44041
$ws.Range("C37").Value = 6513
$ws.Range("D37").Value = 411
$ws.Range("E37").Value = 338
$ws.Range("G37").Value = 5.98
$ws.Range("H37").Value = 2.21
$ws.Range("K37").Value = 5650
$ws.Range("L37").Value = 407

# Row 45 -- Ohio: the scrape failed this run, clear the fetched fields and
# mark the row as an error instead of "Success!"
$ws.Range("B45:H45").Value = ""
$ws.Range("B45:H45").Style = "Normal"
$ws.Range("K45:L45").Value = ""
$ws.Range("K45:L45").Style = "Normal"
$ws.Range("J45").Value = $false
$ws.Range("O45").Value = "An error occurred. ... JSONDecodeError('Expecting value: line 1 column 1 (char 0)')"
